# Base layer rectangles on actual global coord position
# Adds two new columns (N: minR, O: maxR) to Sheet2, filling in the
# per-row rectangle bounds derived from the actual global coordinate
# position of each layer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Header row (row 3) labels for the two new columns.
$ws.Range("N3").Value = "minR"
$ws.Range("O3").Value = "maxR"

# Per-layer minR / maxR values (constant across the 5 modules of each layer).
$minR = @{ 0 = 2945; 1 = 3071; 2 = 3197; 3 = 3323; 4 = 3449; 5 = 3575 }
$maxR = @{ 0 = 3065; 1 = 3191; 2 = 3317; 3 = 3443; 4 = 3569; 5 = 3695 }

$row = 4
foreach ($layer in 0..5) {
    foreach ($module in 0..4) {
        $ws.Cells.Item($row, 14).Value = $minR[$layer]
        $ws.Cells.Item($row, 15).Value = $maxR[$layer]
        $row++
    }
}

# Match the author's final selection state on Sheet2.
$ws.Activate()
$ws.Range("O4:O33").Select()
